# Update "想去人数" (wish-to-go count) values in the "展览" and "全部类型"
# sheets, reflecting a newer data snapshot (F column updates only).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F values to update, keyed by row.
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 1901
    3  = 29
    5  = 111
    8  = 252
    13 = 4497
    16 = 493
    17 = 445
    21 = 2350
    23 = 64
    24 = 43
    26 = 2207
    30 = 158
    31 = 102
    33 = 221
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (all types) - same underlying events, but row indices are
# shifted by 1 starting at row 14 because of an extra row present only here.
$sheetAll = $wb.Worksheets.Item("全部类型")
$allTypesUpdates = @{
    2  = 1901
    3  = 29
    5  = 111
    8  = 252
    14 = 4497
    17 = 493
    18 = 445
    22 = 2350
    24 = 64
    25 = 43
    27 = 2207
    31 = 158
    32 = 102
    34 = 221
}
foreach ($row in $allTypesUpdates.Keys) {
    $sheetAll.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
